$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: switch refArea/dimension labels to iaest-measure labels for comarca/provincia,
# and switch sexo/ccaa/edad-grandes-grupos from iaest-dimension to iaest-measure
$ws.Range("F2").Value = "iaest-measure:residencia-comarca-nombre"
$ws.Range("G2").Value = "iaest-measure:sexo"
$ws.Range("J2").Value = "iaest-measure:residencia-provincia-nombre"
$ws.Range("K2").Value = "iaest-measure:residencia-ccaa-nombre"
$ws.Range("L2").Value = "iaest-measure:edad-grandes-grupos"

# Row 3: these columns become "medida" instead of "dim"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("J3").Value = "medida"
$ws.Range("K3").Value = "medida"
$ws.Range("L3").Value = "medida"

# Row 4: these columns become "xsd:int" instead of skos:Concept/URI-*
$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"
$ws.Range("K4").Value = "xsd:int"
$ws.Range("L4").Value = "xsd:int"

# Row 5: remove mapping file references except for the "ano" column (B5)
$ws.Range("G5").Clear()
$ws.Range("K5").Clear()
$ws.Range("L5").Clear()
